# Add IPC and CPU_util columns (F, G) to the benchmark results table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2) ---
$ws.Range("F2").Value = "IPC "
$ws.Range("G2").Value = "CPU_util"

# --- Data rows (3-11), one pair of (IPC, CPU_util) values per benchmark run ---
$ws.Range("F3").Value = 1.76
$ws.Range("G3").Value = 0.97

$ws.Range("F4").Value = 1.77
$ws.Range("G4").Value = 0.98

$ws.Range("F5").Value = 1.77
$ws.Range("G5").Value = 0.97

$ws.Range("F6").Value = 1.37
$ws.Range("G6").Value = 0.98

$ws.Range("F7").Value = 1.38
$ws.Range("G7").Value = 0.99

$ws.Range("F8").Value = 1.37
$ws.Range("G8").Value = 0.99

$ws.Range("F9").Value = 0.97
$ws.Range("G9").Value = 0.83

$ws.Range("F10").Value = 0.99
$ws.Range("G10").Value = 0.51

$ws.Range("F11").Value = 0.99
$ws.Range("G11").Value = 0.44

# Match the author's final cursor position after adding the new columns.
[void]$ws.Range("G13").Select()

# Best-effort cosmetic follow-up: the source workbook's fonts/style name were
# localized (Chinese "新細明體" / "一般"); the edited copy (opened/saved from a
# different machine) normalizes them to the Calibri / "Normal" defaults.
$wb.Styles.Item(1).Font.Name = "Calibri"
